$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "room_nav_direction" header in K1
$ws.Range("K1").Value = "room_nav_direction"

# Add the new "Room Navigation" row at row 21
$ws.Range("A21").Value = "Room Navigation"
$ws.Range("B21").Value = "Room"
$ws.Range("C21").Value = 125
$ws.Range("D21").Value = 25
$ws.Range("K21").Value = "left"

# Update selection to K1 (matches the saved view state in the diff)
$null = $ws.Range("K1").Select()
